# Add local authority HQ postcode column (D) to the LA website lookup sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> HQ postcode for each local authority (column D).
# Row 1 is the header, row 3 ("Grant aided") has no postcode.
$postcodes = @{
    1  = "postcode"
    2  = "EH1"
    4  = "AB10"
    5  = "AB10"
    6  = "DD8"
    7  = "PA20"
    8  = "FK10"
    9  = "DG1"
    10 = "DD1"
    11 = "KA3"
    12 = "G66"
    13 = "EH41"
    14 = "G46"
    15 = "EH1"
    16 = "HS1"
    17 = "FK1"
    18 = "KY11"
    19 = "G1"
    20 = "AB3"
    21 = "PA15"
    22 = "EH22"
    23 = "AB3"
    24 = "KA12"
    25 = "ML1"
    26 = "KW15"
    27 = "KY13"
    28 = "PA1"
    29 = "EH26"
    30 = "ZE1"
    31 = "KA7"
    32 = "ML3"
    33 = "G63"
    34 = "G82"
    35 = "EH54"
}

foreach ($row in $postcodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $postcodes[$row]
}
